# Atualização de Lista de Requisitos
# Applies the content restructuring described by the diff:
#  - Consolidates "Marcar/Desmarcar/Remarcar Consulta" (UC02-04) into a single
#    "Manter Agenda de Atendimento Médicos" use case (new row 4).
#  - Shifts the remaining use cases up, removing the old UC18-20 rows.
#  - Updates several requirement descriptions.
#  - Removes the now-unused trailing rows (20-22) so the table ends at row 19.
#  - Adjusts row 4's height and the sheet view / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header (rows 1-2) ----------------------------------------------------
$ws.Range("B1").Value = "Requisitos"
$ws.Range("C1").Value = "Caso de Uso"
$ws.Range("E1").Value = "Ator"
$ws.Range("C2").Value = "Nº"
$ws.Range("D2").Value = "Descrição"

# ---- Data rows (3-19) ------------------------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Identificar os usuários: Para que os usuários possam acessar as funcionalidades do sistema, com base em suas respectivas permissões é necessário que estejam `u{201C}logados`u{201D};"
$ws.Range("C3").Value = "UC 01"
$ws.Range("D3").Value = "Logar"
$ws.Range("E3").Value = "Atendente | Gerente | Médico"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Agendar, remarcar, desmarcar e consultar agendamentos de consultas."
$ws.Range("C4").Value = "UC 02"
$ws.Range("D4").Value = "Manter Agenda de Atendimento Médicos."
$ws.Range("E4").Value = "Atendente | Gerente"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "O sistema deverá permitir aos usuários autorizados Cadastrar, Editar, Excluir ou Consultar Pacientes."
$ws.Range("C5").Value = "UC 03"
$ws.Range("D5").Value = "Manter Paciente"
$ws.Range("E5").Value = "Atendente | Gerente"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "O sistema permitirá ao atendente confirmar a presença do paciente e incluí-lo em uma fila de atendimento."
$ws.Range("C6").Value = "UC 04"
$ws.Range("D6").Value = "Registrar Consulta"
$ws.Range("E6").Value = "Atendente | Gerente"

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "O sistema deverá permitir aos usuários autorizados, ao registrar a consulta, gerar, atualizar e/ou consultar o prontuário de cada paciente;"
$ws.Range("C7").Value = "UC 05"
$ws.Range("D7").Value = "Manter Prontuário"
$ws.Range("E7").Value = "Atendente | Gerente | Médico"

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "O sistema deverá permiti aos usuários autorizados, gerar pagamento paras as consultas registradas;"
$ws.Range("C8").Value = "UC 06"
$ws.Range("D8").Value = "Pagar Consulta"
$ws.Range("E8").Value = "Atendente | Gerente"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "O sistema deverá permitir aos usuários autorizados consultar tabela de preços para cada tipo de atendimento;"
$ws.Range("C9").Value = "UC 07"
$ws.Range("D9").Value = "Consultar Tabela de Preços"
$ws.Range("E9").Value = "Atendente | Gerente"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "O sistema deverá permitir aos usuários autorizados solicitar autorização ao Sistema de convênio para Atendimentos pelo o plano de saúde;"
$ws.Range("C10").Value = "UC 08"
$ws.Range("D10").Value = "Autorizar por Convênio"
$ws.Range("E10").Value = "Atendente | Gerente | Convênio"

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "O sistema deverá permitir aos usuários autorizados gerar pagamentos para os atendimentos feitos."
$ws.Range("C11").Value = "UC 09"
$ws.Range("D11").Value = "Gerar Pagamento"
$ws.Range("E11").Value = "Atendente | Gerente"

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "O sistema deverá permitir aos usuários autorizados Cadastrar, Editar, Excluir e Usuários;"
$ws.Range("C12").Value = "UC 10"
$ws.Range("D12").Value = "Manter usuários"
$ws.Range("E12").Value = "Gerente"

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "O sistema deverá permitir aos usuários autorizados Cadastrar, Editar, Excluir e Consultar Tabelas de Preço;"
$ws.Range("C13").Value = "UC 11"
$ws.Range("D13").Value = "Manter Tabela de Preços"
$ws.Range("E13").Value = "Gerente"

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "O sistema deverá permitir aos usuários autorizados Cadastrar, Editar, Excluir e Consultar Convênios;"
$ws.Range("C14").Value = "UC 12"
$ws.Range("D14").Value = "Manter Convênios"
$ws.Range("E14").Value = "Gerente"

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "O sistema deverá permitir aos usuários autorizados Cadastrar, Editar, Excluir e Consultar Horário de Atendimento dos Médicos;"
$ws.Range("C15").Value = "UC 13"
$ws.Range("D15").Value = "Manter Horarios de Atendimentos"
$ws.Range("E15").Value = "Gerente"

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "O sistema deverá permitir aos usuários autorizados Prescrever Medicamento;"
$ws.Range("C16").Value = "UC 14"
$ws.Range("D16").Value = "Prescrever Medicamento"
$ws.Range("E16").Value = "Médico"

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "O sistema deverá permitir aos usuários autorizados Prescrever Exames;"
$ws.Range("C17").Value = "UC 15"
$ws.Range("D17").Value = "Prescrever Exame"
$ws.Range("E17").Value = "Médico"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "O sistema deverá permitir aos usuários autorizados Gerar Atestados Médicos;"
$ws.Range("C18").Value = "UC 16"
$ws.Range("D18").Value = "Gerar Atestado Médico"
$ws.Range("E18").Value = "Médico"

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "O sistema deverá permitir aos usuários autorizados consultar Agendamento de consultas feitas;"
$ws.Range("C19").Value = "UC 17"
$ws.Range("D19").Value = "Consultar Atendimentos agendados"
$ws.Range("E19").Value = "Médico"

# ---- Remove the now-unused rows 20-22 (table shrinks from 20 to 17 items) -
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(20).Delete()

# ---- Row height tweak (row 4 grew a bit to fit the new merged description) -
$ws.Rows.Item(4).RowHeight = 37.5

# ---- View / selection ------------------------------------------------------
$ws.Activate()
$ws.Range("B13").Select()

# ---- Window size metadata (best effort; cosmetic, may not persist) --------
$excel.ActiveWindow.Width = 24000
$excel.ActiveWindow.Height = 9735
